$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 22499.25
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 22499.25
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 22499.25
$ws.Range("N10").Value = -23085.25

# Row 86
$ws.Range("I86").Value = 500003600
$ws.Range("J86").Value = 142858510
$ws.Range("K86").Value = 500003600
$ws.Range("L86").Value = 142858510
$ws.Range("M86").Value = -500002477
$ws.Range("N86").Value = -142860756

# Row 89
$ws.Range("I89").Value = 500003600
$ws.Range("J89").Value = 142858510
$ws.Range("K89").Value = 2500018000
$ws.Range("L89").Value = 714292550
$ws.Range("M89").Value = -2500012384
$ws.Range("N89").Value = -714303782

# Row 103
$ws.Range("H103").Value = 839.03705
$ws.Range("I103").Value = 460.21054
$ws.Range("J103").Value = 1738.75
$ws.Range("K103").Value = 1380.63162
$ws.Range("L103").Value = 5216.25
$ws.Range("M103").Value = -794.6316199999999
$ws.Range("N103").Value = -6388.25

# Row 121
$ws.Range("H121").Value = 2652
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2652
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 7956
$ws.Range("N121").Value = -11450

# Row 132
$ws.Range("H132").Value = 2050.5945
$ws.Range("I132").Value = 1574.1613
$ws.Range("J132").Value = 4512.1665
$ws.Range("K132").Value = 4722.4839
$ws.Range("L132").Value = 13536.4995
$ws.Range("M132").Value = -2192.4839
$ws.Range("N132").Value = -18596.4995

# Row 137
$ws.Range("H137").Value = 50698.824
$ws.Range("I137").Value = 50698.824
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 152096.472
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -149546.472

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 4666.6665
$ws.Range("I22").Value = 4666.6665
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4666.6665
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -4367.6665

# Row 32
$ws.Range("H32").Value = 28760576
$ws.Range("I32").Value = 33965816
$ws.Range("J32").Value = 7939610
$ws.Range("K32").Value = 33965816
$ws.Range("L32").Value = 7939610
$ws.Range("M32").Value = -33965529
$ws.Range("N32").Value = -7940184

# Row 74
$ws.Range("H74").Value = 3001.6155
$ws.Range("I74").Value = 3006.2942
$ws.Range("J74").Value = 2992.7778
$ws.Range("K74").Value = 3006.2942
$ws.Range("L74").Value = 2992.7778
$ws.Range("M74").Value = -2132.2942
$ws.Range("N74").Value = -4740.7778

# Row 77
$ws.Range("H77").Value = 3001.6155
$ws.Range("I77").Value = 3006.2942
$ws.Range("J77").Value = 2992.7778
$ws.Range("K77").Value = 15031.471
$ws.Range("L77").Value = 14963.889
$ws.Range("M77").Value = -10663.471
$ws.Range("N77").Value = -23699.889

# Row 94
$ws.Range("H94").Value = 44166.332
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 44166.332
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 44166.332
$ws.Range("N94").Value = -45968.332

# Row 131
$ws.Range("H131").Value = 117980.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 117980.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 117980.5
$ws.Range("N131").Value = -128060.5

# Row 132
$ws.Range("H132").Value = 4450.3687
$ws.Range("I132").Value = 4283.8667
$ws.Range("J132").Value = 5074.75
$ws.Range("K132").Value = 12851.6001
$ws.Range("L132").Value = 15224.25
$ws.Range("M132").Value = -10321.6001
$ws.Range("N132").Value = -20284.25

$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 242.5
$ws.Range("I11").Value = 242.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 242.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -102.5
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6610.273
$ws.Range("I31").Value = 2657.6
$ws.Range("J31").Value = 8328.825999999999
$ws.Range("K31").Value = 2657.6
$ws.Range("L31").Value = 8328.825999999999
$ws.Range("M31").Value = -2362.6
$ws.Range("N31").Value = -8918.825999999999

# Row 34
$ws.Range("H34").Value = 6610.273
$ws.Range("I34").Value = 2657.6
$ws.Range("J34").Value = 8328.825999999999
$ws.Range("K34").Value = 2657.6
$ws.Range("L34").Value = 8328.825999999999
$ws.Range("M34").Value = -2455.6
$ws.Range("N34").Value = -8732.825999999999

# Row 110
$ws.Range("H110").Value = 101890.664
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 101890.664
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 101890.664
$ws.Range("N110").Value = -110070.664

# Row 132
$ws.Range("H132").Value = 3503.9033
$ws.Range("I132").Value = 3030.1853
$ws.Range("J132").Value = 6701.5
$ws.Range("K132").Value = 9090.555899999999
$ws.Range("L132").Value = 20104.5
$ws.Range("M132").Value = -6560.555899999999
$ws.Range("N132").Value = -25164.5

# Row 140
$ws.Range("H140").Value = 999999
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 999999
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 999999
$ws.Range("N140").Value = -1010359

# Row 141
$ws.Range("H141").Value = 807951
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 807951
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 807951
$ws.Range("N141").Value = -818311

$ws = $wb.Worksheets.Item("CUL")
# Row 133
$ws.Range("H133").Value = 3340.6667
$ws.Range("I133").Value = 2511
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 7533
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -2473
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()

# Row 113
$ws.Range("H113").Value = 26553.389
$ws.Range("I113").Value = 22623.908
$ws.Range("J113").Value = 32728.285
$ws.Range("K113").Value = 22623.908
$ws.Range("L113").Value = 32728.285
$ws.Range("M113").Value = -20453.908
$ws.Range("N113").Value = -37068.285

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8316.333000000001
$ws.Range("I7").Value = 8020.6
$ws.Range("J7").Value = 9795
$ws.Range("K7").Value = 8020.6
$ws.Range("L7").Value = 9795
$ws.Range("M7").Value = -7908.6
$ws.Range("N7").Value = -10019

# Row 16
$ws.Range("H16").Value = 2182.6667
$ws.Range("I16").Value = 1384.4286
$ws.Range("J16").Value = 3300.2
$ws.Range("K16").Value = 1384.4286
$ws.Range("L16").Value = 3300.2
$ws.Range("M16").Value = -1214.4286
$ws.Range("N16").Value = -3640.2

# Row 116
$ws.Range("H116").Value = 158819.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 158819.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 158819.5
$ws.Range("N116").Value = -167997.5

# Row 126
$ws.Range("H126").Value = 8316.333000000001
$ws.Range("I126").Value = 8020.6
$ws.Range("J126").Value = 9795
$ws.Range("K126").Value = 24061.8
$ws.Range("L126").Value = 29385
$ws.Range("M126").Value = -21591.8
$ws.Range("N126").Value = -34325

# Row 132
$ws.Range("H132").Value = 3131.342
$ws.Range("I132").Value = 2874.7188
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 8624.1564
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -6094.1564
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("WVR")
# Row 102
$ws.Range("H102").Value = 84995
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 84995
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 84995
$ws.Range("N102").Value = -91485

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Row 130
$ws.Range("H130").Value = 54544
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 54544
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 54544
$ws.Range("N130").Value = -64584

# Row 132
$ws.Range("H132").Value = 1455.1923
$ws.Range("I132").Value = 1331.9565
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 3995.8695
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -1465.8695
$ws.Range("N132").Value = -12260
